# Refresh the cryptos price/volume table to match the latest scrape.
# Column D holds price strings that often look numeric (e.g. "1.634.67"
# or "1.00"); a leading quote-prefix keeps Excel from re-interpreting them
# as actual numbers (which would silently drop trailing zeros or rewrite
# the text with floating point noise). Column E entries already contain
# padding/percent signs so they are stored as text natively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.969.76"
$ws.Range("D3").Value = "'1.634.67"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'23.49"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'1.865.23"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'1.633.68"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'0.564"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "'65.81"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'27.966.28"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'232.24"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'7.64"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -5.90%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "'155.40"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'15.65"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'1.408.11"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = "  +12.00%  "
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").Value = "'0.557"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'0.868"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'66.90"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.82"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.46"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").Value = "'1.775.87"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'88.26"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "'0.0₆0105"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "'0.0504"
$ws.Range("E51").Value = "  -0.38%  "
